$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "58÷7="
$t.Cell(1, 2).Range.Text = "62÷5="
$t.Cell(1, 3).Range.Text = "80÷6="
$t.Cell(1, 4).Range.Text = "47÷3="
$t.Cell(1, 5).Range.Text = "35÷6="

# Row 5
$t.Cell(5, 1).Range.Text = "36÷9="
$t.Cell(5, 2).Range.Text = "53÷3="
$t.Cell(5, 3).Range.Text = "97÷3="
$t.Cell(5, 4).Range.Text = "16÷6="
$t.Cell(5, 5).Range.Text = "67÷7="

# Row 9
$t.Cell(9, 1).Range.Text = "79÷4="
$t.Cell(9, 2).Range.Text = "28÷7="
$t.Cell(9, 3).Range.Text = "42÷2="
$t.Cell(9, 4).Range.Text = "30÷5="
$t.Cell(9, 5).Range.Text = "56÷5="

# Row 13
$t.Cell(13, 1).Range.Text = "16÷3="
$t.Cell(13, 2).Range.Text = "98÷8="
$t.Cell(13, 3).Range.Text = "59÷2="
$t.Cell(13, 4).Range.Text = "21÷6="
$t.Cell(13, 5).Range.Text = "75÷9="

# Row 17
$t.Cell(17, 1).Range.Text = "49÷9="
$t.Cell(17, 2).Range.Text = "38÷7="
$t.Cell(17, 3).Range.Text = "71÷6="
$t.Cell(17, 4).Range.Text = "50÷5="
$t.Cell(17, 5).Range.Text = "16÷4="
